$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 62 (new row) ---
$ws.Range("C62").Value = "DPDT LEAD FREE SWITCH DC - 6GHZ"
$ws.Range("B62").Value = "4"
$ws.Range("D62").Value = "MASWSS0129TR-3000"
$ws.Range("E62").Value = "1465-1374-1-ND"

# --- Row 58 extra cells ---
$ws.Range("J58").Formula = "=2.6*(0.09+0.04)"
$ws.Range("K58").Value = "watts"

# --- Row 63 (new row) ---
$ws.Range("B63").Value = "1"
$ws.Range("C63").Value = "IC AMP AUDIO PWR 1W MONO 10MSOP"
$ws.Range("D63").Value = "TPA0253DGQR"
$ws.Range("E63").Value = "296-7006-1-ND"

# --- Row 64 (new row) ---
$ws.Range("B64").Value = "2"
$ws.Range("C64").Value = "CAP CER 220UF 6.3V 20% X5R 1210"
$ws.Range("E64").Value = "1276-3375-1-ND"

# --- Row 65 (new row) ---
$ws.Range("B65").Value = "1"
$ws.Range("C65").Value = "IC SWITCH SPDT SC70-6"
$ws.Range("E65").Value = "296-14909-1-ND"
$ws.Range("K65").Formula = "=2.9*0.15"

# --- Row 66 (new row) ---
$ws.Range("B66").Value = "1"
$ws.Range("C66").Value = "IC MULTIPLEXER 2X2 10UMAX"
$ws.Range("E66").Value = "MAX4525CUB+-ND"
$ws.Range("K66").Formula = "=2000/150"

# --- Row 67 (new row) ---
$ws.Range("B67").Value = "1"
$ws.Range("C67").Value = "IC REG BUCK SYNC ADJ 1A SOT25"
$ws.Range("E67").Value = "AP3417CKTR-G1DICT-ND"

# --- Row 68 (new row) ---
$ws.Range("B68").Value = "1"
$ws.Range("C68").Value = "INDUCTOR 2.2UH 4.2A 30% SMD"

# --- back to row 67 for the note in F, then row 68 E ---
$ws.Range("F67").Value = "These are so cheap and small, what's the catch?"
$ws.Range("E68").Value = "587-2098-1-ND"

# --- Row 60/61 numeric additions ---
$ws.Range("K60").Value = 2000
$ws.Range("K61").Formula = "=K60/(90+40)"

# --- Selection state, matching the final saved view ---
$ws.Range("C58").Select()
